$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1 / rId1 -> sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 224
$ws1.Range("F6").Value = 9882
$ws1.Range("F7").Value = 894
$ws1.Range("F10").Value = 3947
$ws1.Range("F16").Value = 558

# Sheet "全部类型" (sheetId=4 / rId4 -> sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 224
$ws4.Range("F7").Value = 9882
$ws4.Range("F8").Value = 894
$ws4.Range("F11").Value = 3947
$ws4.Range("F17").Value = 558
